$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "2024-07-24 13:54:44"
$ws.Range("C8").Value = "2024-07-24 13:55:40"
$ws.Range("C9").Value = "2024-07-24 13:56:03"
$ws.Range("C10").Value = "2024-07-24 13:56:13"
$ws.Range("D10").Value = "2024-07-24 13:56:16"
$ws.Range("E10").Value = "2024-07-24 13:56:22"
